$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

Set-TextCell $ws "D2" "65.469.48"
Set-TextCell $ws "E2" "  -0.31%  "
Set-TextCell $ws "D3" "2.636.99"
Set-TextCell $ws "E3" "  -1.30%  "
Set-TextCell $ws "E4" "  -0.02%  "
Set-TextCell $ws "D5" "593.49"
Set-TextCell $ws "E5" "  -1.02%  "
Set-TextCell $ws "D6" "155.33"
Set-TextCell $ws "E6" "  -0.27%  "
Set-TextCell $ws "E7" "  +0.00%  "
Set-TextCell $ws "D8" "0.623"
Set-TextCell $ws "E8" "  +3.13%  "
Set-TextCell $ws "E9" "  +3.30%  "
Set-TextCell $ws "D10" "0.394"
Set-TextCell $ws "E10" "  -0.38%  "
Set-TextCell $ws "D11" "5.75"
Set-TextCell $ws "E11" "  -2.65%  "
Set-TextCell $ws "E12" "  +0.60%  "
Set-TextCell $ws "D13" "28.57"
Set-TextCell $ws "E13" "  -2.51%  "
Set-TextCell $ws "D14" "0.0000194"
Set-TextCell $ws "E14" "  -0.82%  "
Set-TextCell $ws "D15" "3.110.61"
Set-TextCell $ws "E15" "  -1.28%  "
Set-TextCell $ws "D16" "65.332.09"
Set-TextCell $ws "E16" "  -0.23%  "
Set-TextCell $ws "D17" "2.637.68"
Set-TextCell $ws "E17" "  -1.95%  "
Set-TextCell $ws "D18" "12.50"
Set-TextCell $ws "E18" "  -0.22%  "
Set-TextCell $ws "E19" "  -1.97%  "
Set-TextCell $ws "D20" "7.41"
Set-TextCell $ws "E20" "  -1.16%  "
Set-TextCell $ws "D21" "347.02"
Set-TextCell $ws "E21" "  -0.90%  "
Set-TextCell $ws "E22" "  +0.14%  "
Set-TextCell $ws "D23" "68.74"
Set-TextCell $ws "E23" "  -2.04%  "
Set-TextCell $ws "D24" "0.0000111"
Set-TextCell $ws "E24" "  +2.19%  "
Set-TextCell $ws "D25" "9.55"
Set-TextCell $ws "E25" "  -2.10%  "
Set-TextCell $ws "D26" "1.63"
Set-TextCell $ws "E26" "  -0.28%  "
Set-TextCell $ws "D27" "1.57"
Set-TextCell $ws "E27" "  -2.78%  "
Set-TextCell $ws "D28" "0.163"
Set-TextCell $ws "E28" "  -2.90%  "
Set-TextCell $ws "D29" "1.00"
Set-TextCell $ws "E29" "  -0.01%  "
Set-TextCell $ws "D30" "533.84"
Set-TextCell $ws "E30" "  -0.45%  "
Set-TextCell $ws "B31" "Aptos"
Set-TextCell $ws "C31" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextCell $ws "D31" "7.79"
Set-TextCell $ws "E31" "  -3.62%  "
Set-TextCell $ws "B32" "PancakeSwap"
Set-TextCell $ws "C32" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextCell $ws "D32" "2.10"
Set-TextCell $ws "E32" "  -2.47%  "
Set-TextCell $ws "D33" "1.74"
Set-TextCell $ws "E33" "  -0.68%  "
Set-TextCell $ws "B34" "NEARProtocol"
Set-TextCell $ws "C34" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextCell $ws "D34" "5.42"
Set-TextCell $ws "E34" "  +0.19%  "
Set-TextCell $ws "B35" "RenderToken"
Set-TextCell $ws "C35" "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
Set-TextCell $ws "D35" "6.35"
Set-TextCell $ws "E35" "  -2.72%  "
Set-TextCell $ws "D36" "0.419"
Set-TextCell $ws "E36" "  -0.89%  "
Set-TextCell $ws "D37" "20.22"
Set-TextCell $ws "E37" "  -0.61%  "
Set-TextCell $ws "E38" "  -0.05%  "
Set-TextCell $ws "D39" "1.90"
Set-TextCell $ws "E39" "  -2.11%  "
Set-TextCell $ws "D40" "152.76"
Set-TextCell $ws "E40" "  -4.27%  "
Set-TextCell $ws "E41" "  -0.01%  "
Set-TextCell $ws "D42" "159.32"
Set-TextCell $ws "E42" "  -3.79%  "
Set-TextCell $ws "D43" "4.04"
Set-TextCell $ws "E43" "  -1.05%  "
Set-TextCell $ws "D44" "2.28"
Set-TextCell $ws "E44" "  +2.47%  "
Set-TextCell $ws "D45" "0.0601"
Set-TextCell $ws "E45" "  -1.82%  "
Set-TextCell $ws "D46" "22.39"
Set-TextCell $ws "E46" "  -2.56%  "
Set-TextCell $ws "D47" "0.631"
Set-TextCell $ws "E47" "  -2.51%  "
Set-TextCell $ws "D48" "0.0253"
Set-TextCell $ws "E48" "  -2.46%  "
Set-TextCell $ws "D49" "0.0986"
Set-TextCell $ws "E49" "  -1.09%  "
Set-TextCell $ws "D50" "0.0₆0249"
Set-TextCell $ws "E50" "  +9.13%  "
Set-TextCell $ws "D51" "19.48"
Set-TextCell $ws "E51" "  -2.21%  "
